$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.059.85'
$ws.Range("E2").Value = '  -2.28%  '
$ws.Range("D3").Value = '3.126.99'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.59'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.30%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.119.95'
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.62%  '
$ws.Range("E10").Value = '  -3.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.24'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.28'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.73%  '
$ws.Range("D15").Value = '3.637.11'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("E16").Value = '  +2.03%  '
$ws.Range("D17").Value = '63.038.63'
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("D18").Value = '3.121.15'
$ws.Range("E18").Value = '  -0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.74'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.98'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.28'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.701'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.72'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.12'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.05'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.73'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.17'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.95'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -7.51%  '
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.90'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("E33").Value = '  -8.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.54'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.23%  '
$ws.Range("E35").Value = '  -3.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.01'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").Value = '0.0₃0712'
$ws.Range("E38").Value = '  -5.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0390'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '424.47'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -6.35%  '
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.28'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.70'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -10.98%  '
$ws.Range("D44").Value = '2.886.98'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.267'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.76'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.09%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.113'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.63'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.65%  '
